$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Monte Carlo projection outputs for the auto-enrollment (AE) process.
# Proactive auto-enrollment sweeps additional employees into participation,
# shifting Eligible/Participating counts (and every figure derived from them)
# for projection years 1-5.
# Row 2
$ws.Range("C2").Value = 9403
$ws.Range("D2").Value = 8336
$ws.Range("E2").Value = 0.8865255769435286
$ws.Range("F2").Value = 0.8843624018671759
$ws.Range("G2").Value = 0.09935580614203457
$ws.Range("H2").Value = 0.08786653935921919
$ws.Range("I2").Value = 42296214.91215201
$ws.Range("J2").Value = 14919465.388432
$ws.Range("L2").Value = 14919465.388432
$ws.Range("M2").Value = 57215680.300584
$ws.Range("N2").Value = 799203716.6413001
$ws.Range("O2").Value = 781503909.6373
$ws.Range("P2").Value = 0.01866791292104086
$ws.Range("Q2").Value = 0.0190907111333022

# Row 3
$ws.Range("C3").Value = 9584
$ws.Range("D3").Value = 8489
$ws.Range("E3").Value = 0.8857470784641068
$ws.Range("F3").Value = 0.883075002600645
$ws.Range("G3").Value = 0.09470491223936861
$ws.Range("H3").Value = 0.08363154062207429
$ws.Range("I3").Value = 44258226.35807258
$ws.Range("J3").Value = 16226745.26610873
$ws.Range("L3").Value = 16226745.26610873
$ws.Range("M3").Value = 60484971.62418132
$ws.Range("N3").Value = 844425385.1581769
$ws.Range("O3").Value = 826829949.145153
$ws.Range("P3").Value = 0.01921631626821493
$ws.Range("Q3").Value = 0.01962525097559095

# Row 4
$ws.Range("C4").Value = 9769
$ws.Range("D4").Value = 8666
$ws.Range("E4").Value = 0.8870918210666394
$ws.Range("F4").Value = 0.8839249286005711
$ws.Range("G4").Value = 0.0912970228479114
$ws.Range("H4").Value = 0.0806997144022848
$ws.Range("I4").Value = 46617890.4574882
$ws.Range("J4").Value = 17606478.52025277
$ws.Range("L4").Value = 17606478.52025277
$ws.Range("M4").Value = 64224368.97774097
$ws.Range("N4").Value = 889877250.4101579
$ws.Range("O4").Value = 872313334.6155815
$ws.Range("P4").Value = 0.01978528893972475
$ws.Range("Q4").Value = 0.02018366316503891

# Row 5
$ws.Range("C5").Value = 9967
$ws.Range("D5").Value = 8833
$ws.Range("E5").Value = 0.8862245409852513
$ws.Range("F5").Value = 0.8833883388338833
$ws.Range("G5").Value = 0.08817502547265936
$ws.Range("H5").Value = 0.07789278927892791
$ws.Range("I5").Value = 48849354.23976017
$ws.Range("J5").Value = 18954056.23277476
$ws.Range("L5").Value = 18954056.23277476
$ws.Range("M5").Value = 67803410.47253492
$ws.Range("N5").Value = 936088635.9693686
$ws.Range("O5").Value = 918487413.684114
$ws.Range("P5").Value = 0.02024814264852905
$ws.Range("Q5").Value = 0.0206361632727756

# Row 6
$ws.Range("C6").Value = 10169
$ws.Range("D6").Value = 9041
$ws.Range("E6").Value = 0.8890746386075327
$ws.Range("F6").Value = 0.8865463816434594
$ws.Range("G6").Value = 0.08398739077535673
$ws.Range("H6").Value = 0.07445871739556777
$ws.Range("I6").Value = 50884666.42367596
$ws.Range("J6").Value = 20382863.20657356
$ws.Range("L6").Value = 20382863.20657356
$ws.Range("M6").Value = 71267529.63024952
$ws.Range("N6").Value = 982042596.0921406
$ws.Range("O6").Value = 964335049.8920094
$ws.Range("P6").Value = 0.02075557953156354
$ws.Range("Q6").Value = 0.02113670265210844

